# Deduccion_Planilla_Prueba.xlsx edit
# Removes the "IdInstitucionFinanciera" column (B) and the "idDeducciones"
# column (E), leaving: Identidad, Monto, Comentario, Nombres, Apellidos.
# The surviving "Monto" column (now column B) gets a numeric display
# format instead of the old text format, and the active selection moves
# to the Monto column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rightmost column first (idDeducciones, col E) so column
# indices for the remaining deletion stay stable, then delete
# IdInstitucionFinanciera (col B). This leaves:
#   A Identidad | B Monto | C Comentario | D Nombres | E Apellidos
[void]$ws.Columns.Item(5).Delete()
[void]$ws.Columns.Item(2).Delete()

# The Monto column (now B) should show numbers with 2 decimals instead
# of the generic text format it inherited from its old neighbours.
$ws.Range("B2:B4").NumberFormat = "0.00"

# Match the author's final selection: the Monto column.
[void]$ws.Range("B2:B4").Select()
